$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-24 down to 14-25).
# This makes room for a dedicated "Docentes responsáveis:" value row.
$ws.Rows.Item(13).Insert()

# The newly inserted row 13 comes back completely blank (no formatting).
# Give B13/C13 the same formatting as the B/C columns elsewhere (copy from row 10)
# and then clear out the stray A13 cell that Insert() may have created.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$ws.Application.CutCopyMode = $false

# Row 10 (Objetivos:) previously held the professor's name by mistake;
# replace it with the real (Portuguese) objectives text.
$ws.Range("B10").Value = "Familiarizar os alunos com resultados fundamentais relativos a: integração de funções de uma variável real, cálculo diferencial de funções de n variáveis reais  e suas aplicações."
$ws.Range("C10").Value = "Familiarizar os alunos com resultados fundamentais relativos a: integração de funções de uma variável real, cálculo diferencial de funções de n variáveis reais  e suas aplicações."

# Row 13 (new, under "Docentes responsáveis:") gets the professor's name.
$ws.Range("B13").Value = "8822123 - Roberta Veloso Garcia"
$ws.Range("C13").Value = "8822123 - Roberta Veloso Garcia"

# Row 14 (Programa resumido:) previously just said "Semestral"; fill with the
# actual (Portuguese) short-syllabus summary.
$ws.Range("B14").Value = "Integração de funções de uma variável real. Funções reais de variáveis reais, Diferenciabilidade, Derivada direcional. Máximos e Mínios em domínios abertos e Multiplicadores de Lagrange"
$ws.Range("C14").Value = "Integração de funções de uma variável real. Funções reais de variáveis reais, Diferenciabilidade, Derivada direcional. Máximos e Mínios em domínios abertos e Multiplicadores de Lagrange"

# Row 16 (Programa:) previously duplicated the "Ativação" date by mistake;
# replace it with the full (Portuguese) syllabus text.
$ws.Range("B16").Value = "Integração de funções reais: Primitivas (Integral indefinida), Integral de Riemann (Integral definida), Teorema fundamental do cálculo, Técnicas de integração e aplicações. O espaço euclidiano R^n: Conjuntos abertos, fechados e compactos.Funções de n várias variáveis Reais: Gráficos e curvas de nível de funções de duas variáveis.Limites e Continuidade: Teorema de WeierstrassDiferenciabilidade: Derivadas parciais, diferencial total, derivadas parciais de ordem superior, teorema de Schwarz, regra da cadeia, planos tangentes e aproximações lineares, derivada direcional, vetor gradiente, teorema da função implícita, jacobiano.Máximos e mínimos: Valores Extremos de funções de duas ou mais variáveis em domínios abertos, Hessiano de uma função real de n variáveis, multiplicadores de Lagrange."
$ws.Range("C16").Value = "Integração de funções reais: Primitivas (Integral indefinida), Integral de Riemann (Integral definida), Teorema fundamental do cálculo, Técnicas de integração e aplicações. O espaço euclidiano R^n: Conjuntos abertos, fechados e compactos.Funções de n várias variáveis Reais: Gráficos e curvas de nível de funções de duas variáveis.Limites e Continuidade: Teorema de WeierstrassDiferenciabilidade: Derivadas parciais, diferencial total, derivadas parciais de ordem superior, teorema de Schwarz, regra da cadeia, planos tangentes e aproximações lineares, derivada direcional, vetor gradiente, teorema da função implícita, jacobiano.Máximos e mínimos: Valores Extremos de funções de duas ou mais variáveis em domínios abertos, Hessiano de uma função real de n variáveis, multiplicadores de Lagrange."

# Row 19 (Método:) previously duplicated the professor's name by mistake;
# replace it with the actual evaluation-method text.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 (Critério:) previously held the evaluation-method text; replace with
# the pass criterion text.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21 (Norma de recuperação:) previously held the pass-criterion text;
# replace with the make-up (recuperação) norm text.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# Row 22 (Bibliografia:) previously held the recovery-norm text; replace with
# the actual bibliography list.
$ws.Range("B22").Value = "GUIDORIZZI, Hamilton L. UM CURSO DE CÁLCULO,2011, 5. ed., v.2LEITHOLD, Louis. CÁLCULO COM GEOMETRIA ANALÍTICA, São Paulo: HARBRA LTDA, 1990. v.2ANTON, Howard; BIVENS, Irl, DAVIS, Stephen. CÁLCULO, 8. ed. São Paulo:Pearson, 2011, v.2SIMMONS, George F. CÁLCULO COM GEOMETRIA ANALÍTICA, São Paulo: Pearson, 2014. v.2STEWART, James. CÁLCULO. revisão técnica Ricardo Miranda Martins. 7. ed. São Paulo: Cengage Learning, 2013. v.2THOMAS, George B. WEIR, Maurice D.; HASS, Joel; GIORDANO, CÁLCULO. revisão técnica Cláudio Hirofume Asano .12.ed. São Paulo: Pearson Education do Brasil, 2013. v.2"
$ws.Range("C22").Value = "GUIDORIZZI, Hamilton L. UM CURSO DE CÁLCULO,2011, 5. ed., v.2LEITHOLD, Louis. CÁLCULO COM GEOMETRIA ANALÍTICA, São Paulo: HARBRA LTDA, 1990. v.2ANTON, Howard; BIVENS, Irl, DAVIS, Stephen. CÁLCULO, 8. ed. São Paulo:Pearson, 2011, v.2SIMMONS, George F. CÁLCULO COM GEOMETRIA ANALÍTICA, São Paulo: Pearson, 2014. v.2STEWART, James. CÁLCULO. revisão técnica Ricardo Miranda Martins. 7. ed. São Paulo: Cengage Learning, 2013. v.2THOMAS, George B. WEIR, Maurice D.; HASS, Joel; GIORDANO, CÁLCULO. revisão técnica Cláudio Hirofume Asano .12.ed. São Paulo: Pearson Education do Brasil, 2013. v.2"

# Row heights: the inserted row 13 already comes back at the default height
# (no custom height flag), matching the rest of the small single-line rows,
# so nothing needs to be done for it.

# Rows that used to be "short" (no explicit custom height) but now hold long
# paragraph text need the larger custom heights seen in similar long rows.
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 120
